$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: Volume/Number and week-of dates ---
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# --- Crime Complaints data table updates (rows 15-30) ---

$ws.Range("F15").Value = "'0"
$ws.Range("N15").Value = -75

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = -63.636363636363
$ws.Range("N16").Value = -89.743589743589

$ws.Range("D17").Value = 4
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -85.714285714285
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = -60
$ws.Range("M17").Value = -27.272727272727
$ws.Range("N17").Value = -61.904761904761

$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "***.*"
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 0
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -55.555555555555
$ws.Range("N18").Value = -94.594594594594

$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 13
$ws.Range("J19").Value = 21
$ws.Range("K19").Value = -38.095238095238
$ws.Range("L19").Value = -51.851851851851
$ws.Range("M19").Value = -18.75
$ws.Range("N19").Value = -27.777777777777

$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -90
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = -85.714285714285
$ws.Range("N20").Value = -92.857142857142

$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = -80
$ws.Range("G21").Value = 39
$ws.Range("H21").Value = -58.974358974359
$ws.Range("I21").Value = 32
$ws.Range("J21").Value = 64
$ws.Range("K21").Value = -50
$ws.Range("L21").Value = -44.827586206896
$ws.Range("M21").Value = -36
$ws.Range("N21").Value = -82.795698924731

$ws.Range("G23").Value = 3
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = -50
$ws.Range("M23").Value = -71.428571428571

$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -27.272727272727
$ws.Range("F24").Value = 28
$ws.Range("G24").Value = 25
$ws.Range("H24").Value = 12
$ws.Range("I24").Value = 54
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 8
$ws.Range("L24").Value = -32.5
$ws.Range("M24").Value = 25.581395348837

$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 25
$ws.Range("I25").Value = 29
$ws.Range("J25").Value = 22
$ws.Range("K25").Value = 31.818181818181
$ws.Range("L25").Value = -35.555555555555

$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = -21.428571428571
$ws.Range("I26").Value = 24
$ws.Range("J26").Value = 28
$ws.Range("K26").Value = -14.285714285714
$ws.Range("L26").Value = -11.111111111111
$ws.Range("M26").Value = -36.842105263157

$ws.Range("F27").Value = "'0"

$ws.Range("F28").Value = "'0"
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -100

$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "***.*"

$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "***.*"

Write-Host "Applied weekly crime data update"
